$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update row 2 "Accept" value (IT ES -> IT US)
$ws.Range("F2").Value = "IT US"

# 2. Insert a new row at position 3 (shifts old rows 3-8 down to 4-9),
#    duplicating the HTTP/s1.rivetweb.org service entry but scoped to the
#    pastorious target host (more efficient ACL grouping).
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "HTTP"
$ws.Range("B3").Value = 443
$ws.Range("C3").Value = "s1.rivetweb.org"
$ws.Range("D3").Value = "shorter.rivetweb.org"
$ws.Range("E3").Value = 80
$ws.Range("F3").Value = "pastorious.rivetweb.org"
$ws.Range("G3").Value = "ALL"
$ws.Range("H3").Value = "enable"

# 3. Row 7 (the s2.rivetweb.org / pastorious HTTP rule, shifted down from
#    the old row 6) now accepts from ALL and rejects a specific flagged IP
#    instead of filtering by country code.
$ws.Range("F7").Value = "ALL"
$ws.Range("G7").Value = "167.172.169.214 RU"

# 4. Widen the "Reject" column to fit the new longer values.
$ws.Columns.Item(7).ColumnWidth = 27.6867

# 5. Restore the active selection to match the edited cell.
$ws.Range("G7").Select()
